$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 787
$ws.Range("I8").Value = 787
$ws.Range("K8").Value = 2361
$ws.Range("M8").Value = -2222
$ws.Range("H40").Value = 2057.1428
$ws.Range("I40").Value = 2069.2307
$ws.Range("J40").Value = 1900
$ws.Range("K40").Value = 2069.2307
$ws.Range("L40").Value = 1900
$ws.Range("M40").Value = -1894.2307
$ws.Range("N40").Value = -2250
$ws.Range("H64").Value = 4051.7258
$ws.Range("I64").Value = 3426.6978
$ws.Range("J64").Value = 5466.263
$ws.Range("K64").Value = 3426.6978
$ws.Range("L64").Value = 5466.263
$ws.Range("M64").Value = -3178.6978
$ws.Range("N64").Value = -5962.263
$ws.Range("H67").Value = 4051.7258
$ws.Range("I67").Value = 3426.6978
$ws.Range("J67").Value = 5466.263
$ws.Range("K67").Value = 3426.6978
$ws.Range("L67").Value = 5466.263
$ws.Range("M67").Value = -2568.6978
$ws.Range("N67").Value = -7182.263
$ws.Range("H74").Value = 3772.4285
$ws.Range("I74").Value = 3236.6365
$ws.Range("J74").Value = 4361.8
$ws.Range("K74").Value = 3236.6365
$ws.Range("L74").Value = 4361.8
$ws.Range("M74").Value = -2300.6365
$ws.Range("N74").Value = -6233.8
$ws.Range("H76").Value = 3614.975
$ws.Range("I76").Value = 3034.6775
$ws.Range("J76").Value = 5613.778
$ws.Range("K76").Value = 3034.6775
$ws.Range("L76").Value = 5613.778
$ws.Range("M76").Value = -2719.6775
$ws.Range("N76").Value = -6243.778
$ws.Range("H77").Value = 3772.4285
$ws.Range("I77").Value = 3236.6365
$ws.Range("J77").Value = 4361.8
$ws.Range("K77").Value = 16183.1825
$ws.Range("L77").Value = 21809
$ws.Range("M77").Value = -11503.1825
$ws.Range("N77").Value = -31169
$ws.Range("H79").Value = 3614.975
$ws.Range("I79").Value = 3034.6775
$ws.Range("J79").Value = 5613.778
$ws.Range("K79").Value = 3034.6775
$ws.Range("L79").Value = 5613.778
$ws.Range("M79").Value = -1942.6775
$ws.Range("N79").Value = -7797.778
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H61").Value = 12003.5
$ws.Range("I61").Value = 6500
$ws.Range("J61").Value = 17507
$ws.Range("K61").Value = 6500
$ws.Range("L61").Value = 17507
$ws.Range("M61").Value = -6288
$ws.Range("N61").Value = -17931
$ws.Range("H63").Value = 8586.666999999999
$ws.Range("I63").Value = 8985.714
$ws.Range("K63").Value = 8985.714
$ws.Range("M63").Value = -8299.714
$ws.Range("H66").Value = 8586.666999999999
$ws.Range("I66").Value = 8985.714
$ws.Range("K66").Value = 44928.57
$ws.Range("M66").Value = -41496.57
$ws.Range("H74").Value = 1944.8
$ws.Range("I74").Value = 3264.8
$ws.Range("J74").Value = 1504.8
$ws.Range("K74").Value = 3264.8
$ws.Range("L74").Value = 1504.8
$ws.Range("M74").Value = -2390.8
$ws.Range("N74").Value = -3252.8
$ws.Range("H77").Value = 1944.8
$ws.Range("I77").Value = 3264.8
$ws.Range("J77").Value = 1504.8
$ws.Range("K77").Value = 16324
$ws.Range("L77").Value = 7524
$ws.Range("M77").Value = -11956
$ws.Range("N77").Value = -16260
$ws.Range("H92").Value = 27933.334
$ws.Range("J92").Value = 27933.334
$ws.Range("L92").Value = 27933.334
$ws.Range("N92").Value = -32925.334
$ws.Range("H102").Value = 1983.9412
$ws.Range("I102").Value = 2013.8
$ws.Range("J102").Value = 1941.2858
$ws.Range("K102").Value = 2013.8
$ws.Range("L102").Value = 1941.2858
$ws.Range("M102").Value = -391.8
$ws.Range("N102").Value = -5185.2858
$ws.Range("H136").Value = 12003.5
$ws.Range("I136").Value = 6500
$ws.Range("J136").Value = 17507
$ws.Range("K136").Value = 19500
$ws.Range("L136").Value = 52521
$ws.Range("M136").Value = -16950
$ws.Range("N136").Value = -57621
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1281.2693
$ws.Range("I99").Value = 898.6667
$ws.Range("J99").Value = 1803
$ws.Range("K99").Value = 898.6667
$ws.Range("L99").Value = 1803
$ws.Range("M99").Value = 599.3333
$ws.Range("N99").Value = -4799
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1736.4445
$ws.Range("I31").Value = 1452.9375
$ws.Range("J31").Value = 4004.5
$ws.Range("K31").Value = 1452.9375
$ws.Range("L31").Value = 4004.5
$ws.Range("M31").Value = -1157.9375
$ws.Range("N31").Value = -4594.5
$ws.Range("H34").Value = 1736.4445
$ws.Range("I34").Value = 1452.9375
$ws.Range("J34").Value = 4004.5
$ws.Range("K34").Value = 1452.9375
$ws.Range("L34").Value = 4004.5
$ws.Range("M34").Value = -1250.9375
$ws.Range("N34").Value = -4408.5
$ws.Range("H62").Value = 4301.6665
$ws.Range("I62").Value = 6202.5
$ws.Range("J62").Value = 3351.25
$ws.Range("K62").Value = 6202.5
$ws.Range("L62").Value = 3351.25
$ws.Range("M62").Value = -5578.5
$ws.Range("N62").Value = -4599.25
$ws.Range("H65").Value = 4301.6665
$ws.Range("I65").Value = 6202.5
$ws.Range("J65").Value = 3351.25
$ws.Range("K65").Value = 31012.5
$ws.Range("L65").Value = 16756.25
$ws.Range("M65").Value = -27892.5
$ws.Range("N65").Value = -22996.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 27936.158
$ws.Range("J9").Value = 87481.164
$ws.Range("L9").Value = 262443.492
$ws.Range("N9").Value = -262891.492
$ws.Range("H56").Value = 6836.6665
$ws.Range("I56").Value = 6836.6665
$ws.Range("K56").Value = 6836.6665
$ws.Range("M56").Value = -6306.6665
$ws.Range("H131").Value = 2617.3088
$ws.Range("I131").Value = 649
$ws.Range("J131").Value = 2708.1538
$ws.Range("K131").Value = 1947
$ws.Range("L131").Value = 8124.4614
$ws.Range("M131").Value = 3093
$ws.Range("N131").Value = -18204.4614
$ws.Range("H132").Value = 1539.9412
$ws.Range("I132").Value = 1766.6666
$ws.Range("J132").Value = 1284.875
$ws.Range("K132").Value = 15899.9994
$ws.Range("L132").Value = 11563.875
$ws.Range("M132").Value = -13369.9994
$ws.Range("N132").Value = -16623.875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5945094
$ws.Range("I14").Value = 7307965
$ws.Range("J14").Value = 39319.668
$ws.Range("K14").Value = 7307965
$ws.Range("L14").Value = 39319.668
$ws.Range("M14").Value = -7307797
$ws.Range("N14").Value = -39655.668
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H70").Value = 26040.043
$ws.Range("I70").Value = 49456
$ws.Range("J70").Value = 4575.4165
$ws.Range("K70").Value = 49456
$ws.Range("L70").Value = 4575.4165
$ws.Range("M70").Value = -49186
$ws.Range("N70").Value = -5115.4165
$ws.Range("H73").Value = 26040.043
$ws.Range("I73").Value = 49456
$ws.Range("J73").Value = 4575.4165
$ws.Range("K73").Value = 49456
$ws.Range("L73").Value = 4575.4165
$ws.Range("M73").Value = -48520
$ws.Range("N73").Value = -6447.4165
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 36700
$ws.Range("J74").Value = 36700
$ws.Range("L74").Value = 36700
$ws.Range("N74").Value = -38696
$ws.Range("H77").Value = 36700
$ws.Range("J77").Value = 36700
$ws.Range("L77").Value = 110100
$ws.Range("N77").Value = -120084
$ws.Range("H93").Value = 7833.357
$ws.Range("I93").Value = 8397.462
$ws.Range("J93").Value = 500
$ws.Range("K93").Value = 8397.462
$ws.Range("L93").Value = 500
$ws.Range("M93").Value = -7149.462
$ws.Range("N93").Value = -2996
$ws.Range("H100").Value = 1468.75
$ws.Range("I100").Value = 1468.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1468.75
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -927.75
$ws.Range("N100").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 2099.6667
$ws.Range("I13").Value = 299
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 299
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -159
$ws.Range("N13").Value = -3280
$ws.Range("H136").Value = 6673.826
$ws.Range("I136").Value = 1154.2222
$ws.Range("J136").Value = 10222.143
$ws.Range("K136").Value = 3462.6666
$ws.Range("L136").Value = 30666.429
$ws.Range("M136").Value = -912.6665999999996
$ws.Range("N136").Value = -35766.429
